$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.308.94"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "3.399.65"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +11.50%  "
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("E12").Value = "  +5.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "685.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").Value = "3.946.77"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "69.357.22"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "3.378.04"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.911"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "557.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.52%  "
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "3.653.54"
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("E38").Value = "  +5.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "0.0₃0727"
$ws.Range("E40").Value = "  +8.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0428"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.338"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
